$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.596199999999991
$ws.Range("B4").Value = 4.805400000000002
$ws.Range("C6").Value = -11.0346
$ws.Range("B7").Value = 6.131400000000002
$ws.Range("C7").Value = -11.0925
$ws.Range("B8").Value = 5.187099999999998
$ws.Range("C8").Value = -10.64609999999999
$ws.Range("A11").Value = -22.17160000000003
$ws.Range("A12").Value = -22.6291
$ws.Range("B12").Value = 5.739700000000004
$ws.Range("B14").Value = 9.447800000000004
$ws.Range("A15").Value = -21.53550000000002
$ws.Range("C19").Value = -12.89579999999999
$ws.Range("C21").Value = -12.6816
$ws.Range("B22").Value = 5.296100000000004
$ws.Range("C24").Value = -11.4735
$ws.Range("C25").Value = -10.72779999999999
